$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 46202.332
$ws.Range("J3").Value = 46202.332
$ws.Range("L3").Value = 46202.332
$ws.Range("N3").Value = -46430.332
$ws.Range("H4").Value = 106.38461
$ws.Range("J4").Value = 299.5
$ws.Range("L4").Value = 299.5
$ws.Range("N4").Value = -527.5
$ws.Range("H7").Value = 10005
$ws.Range("I7").Value = 10005
$ws.Range("K7").Value = 10005
$ws.Range("M7").Value = -9893
$ws.Range("H13").Value = 8997.5
$ws.Range("I13").Value = 7999.5
$ws.Range("K13").Value = 7999.5
$ws.Range("M13").Value = -7830.5
$ws.Range("H14").Value = 10005
$ws.Range("I14").Value = 10005
$ws.Range("K14").Value = 10005
$ws.Range("M14").Value = -9814
$ws.Range("H19").Value = 948.7857
$ws.Range("J19").Value = 1045
$ws.Range("L19").Value = 1045
$ws.Range("N19").Value = -1395
$ws.Range("H33").Value = 1048.8823
$ws.Range("I33").Value = 200.25
$ws.Range("J33").Value = 1803.2222
$ws.Range("K33").Value = 200.25
$ws.Range("L33").Value = 1803.2222
$ws.Range("M33").Value = 28.75
$ws.Range("N33").Value = -2261.2222
$ws.Range("H40").Value = 4484.7144
$ws.Range("I40").Value = 2999
$ws.Range("J40").Value = 4732.3335
$ws.Range("K40").Value = 2999
$ws.Range("L40").Value = 4732.3335
$ws.Range("M40").Value = -2824
$ws.Range("N40").Value = -5082.3335
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10924
$ws.Range("H58").Value = 22732718
$ws.Range("I58").Value = 27778342
$ws.Range("J58").Value = 27408
$ws.Range("K58").Value = 83335026
$ws.Range("L58").Value = 82224
$ws.Range("M58").Value = -83334876
$ws.Range("N58").Value = -82524
$ws.Range("H62").Value = 34505.242
$ws.Range("I62").Value = 118885.43
$ws.Range("K62").Value = 118885.43
$ws.Range("M62").Value = -118261.43
$ws.Range("H65").Value = 34505.242
$ws.Range("I65").Value = 118885.43
$ws.Range("K65").Value = 594427.1499999999
$ws.Range("M65").Value = -591307.1499999999
$ws.Range("H98").Value = 7765.909
$ws.Range("I98").Value = 7765.909
$ws.Range("K98").Value = 7765.909
$ws.Range("M98").Value = -6267.909
$ws.Range("H102").Value = 46202.332
$ws.Range("J102").Value = 46202.332
$ws.Range("L102").Value = 46202.332
$ws.Range("N102").Value = -52692.332
$ws.Range("H106").Value = 2805.6191
$ws.Range("I106").Value = 2620.9
$ws.Range("K106").Value = 2620.9
$ws.Range("M106").Value = -1989.9
$ws.Range("H116").Value = 3753.1428
$ws.Range("I116").Value = 3314.8
$ws.Range("K116").Value = 3314.8
$ws.Range("M116").Value = 127.1999999999998
$ws.Range("H122").Value = 7765.909
$ws.Range("I122").Value = 7765.909
$ws.Range("K122").Value = 23297.727
$ws.Range("M122").Value = -20847.727
$ws.Range("H129").Value = 1739.4
$ws.Range("I129").Value = 982.3333
$ws.Range("K129").Value = 2946.9999
$ws.Range("M129").Value = 2053.0001
$ws.Range("H131").Value = 7584.875
$ws.Range("I131").Value = 6736.3
$ws.Range("K131").Value = 20208.9
$ws.Range("M131").Value = -15168.9
$ws.Range("H132").Value = 1448
$ws.Range("I132").Value = 1343.0731
$ws.Range("J132").Value = 2165
$ws.Range("K132").Value = 4029.2193
$ws.Range("L132").Value = 6495
$ws.Range("M132").Value = -1499.2193
$ws.Range("N132").Value = -11555
$ws.Range("H137").Value = 1731.9231
$ws.Range("I137").Value = 1731.9231
$ws.Range("K137").Value = 5195.7693
$ws.Range("M137").Value = -2645.7693
$ws.Range("H138").Value = 2625.7048
$ws.Range("J138").Value = 3406.4285
$ws.Range("L138").Value = 10219.2855
$ws.Range("N138").Value = -20499.2855
$ws.Range("H141").Value = 2740.5715
$ws.Range("I141").Value = 2536.8
$ws.Range("J141").Value = 3250
$ws.Range("K141").Value = 7610.400000000001
$ws.Range("L141").Value = 9750
$ws.Range("M141").Value = -2430.400000000001
$ws.Range("N141").Value = -20110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1563.5186
$ws.Range("I2").Value = 827.1579
$ws.Range("K2").Value = 827.1579
$ws.Range("M2").Value = -714.1579
$ws.Range("H32").Value = 2640.291
$ws.Range("I32").Value = 2494.9412
$ws.Range("K32").Value = 2494.9412
$ws.Range("M32").Value = -2207.9412
$ws.Range("H33").Value = 52111.11
$ws.Range("I33").Value = 56333.332
$ws.Range("J33").Value = 50000
$ws.Range("K33").Value = 56333.332
$ws.Range("L33").Value = 50000
$ws.Range("M33").Value = -56004.332
$ws.Range("N33").Value = -50658
$ws.Range("H36").Value = 4314.4287
$ws.Range("I36").Value = 4616.8335
$ws.Range("J36").Value = 2500
$ws.Range("K36").Value = 4616.8335
$ws.Range("L36").Value = 2500
$ws.Range("M36").Value = -4270.8335
$ws.Range("N36").Value = -3192
$ws.Range("H45").Value = 1999.2222
$ws.Range("I45").Value = 1998.25
$ws.Range("K45").Value = 1998.25
$ws.Range("M45").Value = -1621.25
$ws.Range("H61").Value = 4688.636
$ws.Range("I61").Value = 4688.636
$ws.Range("K61").Value = 4688.636
$ws.Range("M61").Value = -4476.636
$ws.Range("H97").Value = 1359.3529
$ws.Range("I97").Value = 1240.6
$ws.Range("K97").Value = 1240.6
$ws.Range("M97").Value = -744.5999999999999
$ws.Range("H116").Value = 1563.5186
$ws.Range("I116").Value = 827.1579
$ws.Range("K116").Value = 827.1579
$ws.Range("M116").Value = 1466.8421
$ws.Range("H122").Value = 2198.1738
$ws.Range("I122").Value = 2037.1111
$ws.Range("K122").Value = 6111.3333
$ws.Range("M122").Value = -3661.3333
$ws.Range("H132").Value = 3375
$ws.Range("I132").Value = 3375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7595
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 40166.5
$ws.Range("J134").Value = 40166.5
$ws.Range("L134").Value = 40166.5
$ws.Range("N134").Value = -50306.5
$ws.Range("H136").Value = 4688.636
$ws.Range("I136").Value = 4688.636
$ws.Range("K136").Value = 14065.908
$ws.Range("M136").Value = -11515.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1563.5186
$ws.Range("I3").Value = 827.1579
$ws.Range("K3").Value = 827.1579
$ws.Range("M3").Value = -713.1579
$ws.Range("H17").Value = 10166.667
$ws.Range("I17").Value = 10000
$ws.Range("J17").Value = 10200
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 10200
$ws.Range("M17").Value = -9828
$ws.Range("N17").Value = -10544
$ws.Range("H86").Value = 1737.6666
$ws.Range("I86").Value = 1769.1111
$ws.Range("K86").Value = 1769.1111
$ws.Range("M86").Value = -646.1111000000001
$ws.Range("H89").Value = 1737.6666
$ws.Range("I89").Value = 1769.1111
$ws.Range("K89").Value = 8845.5555
$ws.Range("M89").Value = -3229.5555
$ws.Range("H94").Value = 3846.9473
$ws.Range("I94").Value = 4045.75
$ws.Range("J94").Value = 3506.1428
$ws.Range("K94").Value = 4045.75
$ws.Range("L94").Value = 3506.1428
$ws.Range("M94").Value = -3594.75
$ws.Range("N94").Value = -4408.1428
$ws.Range("H132").Value = 74998.5
$ws.Range("J132").Value = 74998.5
$ws.Range("L132").Value = 74998.5
$ws.Range("N132").Value = -85118.5
$ws.Range("H134").Value = 2062.9412
$ws.Range("I134").Value = 1900.9286
$ws.Range("J134").Value = 2819
$ws.Range("K134").Value = 5702.7858
$ws.Range("L134").Value = 8457
$ws.Range("M134").Value = -3167.7858
$ws.Range("N134").Value = -13527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 2029.8
$ws.Range("J8").Value = 2029.8
$ws.Range("L8").Value = 2029.8
$ws.Range("N8").Value = -2309.8
$ws.Range("H15").Value = 600
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 600
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -940
$ws.Range("H31").Value = 6470.52
$ws.Range("I31").Value = 3166.2666
$ws.Range("J31").Value = 11426.9
$ws.Range("K31").Value = 3166.2666
$ws.Range("L31").Value = 11426.9
$ws.Range("M31").Value = -2871.2666
$ws.Range("N31").Value = -12016.9
$ws.Range("H32").Value = 309.5
$ws.Range("I32").Value = 309.5
$ws.Range("K32").Value = 309.5
$ws.Range("M32").Value = 6.5
$ws.Range("H34").Value = 6470.52
$ws.Range("I34").Value = 3166.2666
$ws.Range("J34").Value = 11426.9
$ws.Range("K34").Value = 3166.2666
$ws.Range("L34").Value = 11426.9
$ws.Range("M34").Value = -2964.2666
$ws.Range("N34").Value = -11830.9
$ws.Range("H41").Value = 27752.715
$ws.Range("I41").Value = 4111.8
$ws.Range("K41").Value = 4111.8
$ws.Range("M41").Value = -3683.8
$ws.Range("H51").Value = 24428.215
$ws.Range("J51").Value = 24428.215
$ws.Range("L51").Value = 24428.215
$ws.Range("N51").Value = -25900.215
$ws.Range("H52").Value = 100172.25
$ws.Range("J52").Value = 99993.336
$ws.Range("L52").Value = 99993.336
$ws.Range("N52").Value = -100581.336
$ws.Range("H55").Value = 30100
$ws.Range("I55").Value = 33020
$ws.Range("J55").Value = 15500
$ws.Range("K55").Value = 33020
$ws.Range("L55").Value = 15500
$ws.Range("M55").Value = -32705
$ws.Range("N55").Value = -16130
$ws.Range("H56").Value = 15000
$ws.Range("I56").Value = 15000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -14155
$ws.Range("N56").ClearContents()
$ws.Range("H58").Value = 2224.4614
$ws.Range("I58").Value = 2224.4614
$ws.Range("K58").Value = 2224.4614
$ws.Range("M58").Value = -2021.4614
$ws.Range("H61").Value = 24428.215
$ws.Range("J61").Value = 24428.215
$ws.Range("L61").Value = 24428.215
$ws.Range("N61").Value = -25124.215
$ws.Range("H69").Value = 18599.8
$ws.Range("I69").Value = 12249.5
$ws.Range("K69").Value = 12249.5
$ws.Range("M69").Value = -11500.5
$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50630
$ws.Range("H72").Value = 18599.8
$ws.Range("I72").Value = 12249.5
$ws.Range("K72").Value = 36748.5
$ws.Range("M72").Value = -33004.5
$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -52184
$ws.Range("H75").Value = 19320
$ws.Range("J75").Value = 19320
$ws.Range("L75").Value = 19320
$ws.Range("N75").Value = -21316
$ws.Range("H78").Value = 19320
$ws.Range("J78").Value = 19320
$ws.Range("L78").Value = 57960
$ws.Range("N78").Value = -67944
$ws.Range("H103").Value = 11670.5
$ws.Range("I103").Value = 11670.5
$ws.Range("K103").Value = 11670.5
$ws.Range("M103").Value = -10498.5
$ws.Range("H132").Value = 2701.606
$ws.Range("I132").Value = 2756.5483
$ws.Range("K132").Value = 8269.644899999999
$ws.Range("M132").Value = -5739.644899999999
$ws.Range("H134").Value = 2407.1428
$ws.Range("I134").Value = 2421.8064
$ws.Range("K134").Value = 7265.4192
$ws.Range("M134").Value = -4730.4192
$ws.Range("H136").Value = 2224.4614
$ws.Range("I136").Value = 2224.4614
$ws.Range("K136").Value = 6673.3842
$ws.Range("M136").Value = -4123.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 128.4
$ws.Range("J2").Value = 55.57143
$ws.Range("L2").Value = 333.42858
$ws.Range("N2").Value = -559.42858
$ws.Range("H16").Value = 449.75
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1327
$ws.Range("H55").Value = 838550
$ws.Range("J55").Value = 6260
$ws.Range("L55").Value = 18780
$ws.Range("N55").Value = -19134
$ws.Range("H107").Value = 2055.6924
$ws.Range("I107").Value = 3181.6667
$ws.Range("J107").Value = 1090.5714
$ws.Range("K107").Value = 9545.000100000001
$ws.Range("L107").Value = 3271.7142
$ws.Range("M107").Value = -7625.000100000001
$ws.Range("N107").Value = -7111.7142
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H129").Value = 3847
$ws.Range("I129").Value = 1092.25
$ws.Range("K129").Value = 3276.75
$ws.Range("M129").Value = 1723.25
$ws.Range("H131").Value = 2358.9375
$ws.Range("I131").Value = 1612
$ws.Range("J131").Value = 2531.3076
$ws.Range("K131").Value = 4836
$ws.Range("L131").Value = 7593.9228
$ws.Range("M131").Value = 204
$ws.Range("N131").Value = -17673.9228
$ws.Range("H137").Value = 5544.4165
$ws.Range("I137").Value = 1761.3
$ws.Range("J137").Value = 6999.4614
$ws.Range("K137").Value = 5283.9
$ws.Range("L137").Value = 20998.3842
$ws.Range("M137").Value = -183.8999999999996
$ws.Range("N137").Value = -31198.3842
$ws.Range("H140").Value = 1648.0769
$ws.Range("I140").Value = 1311.3636
$ws.Range("J140").Value = 3500
$ws.Range("K140").Value = 3934.0908
$ws.Range("L140").Value = 10500
$ws.Range("M140").Value = 1245.9092
$ws.Range("N140").Value = -20860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 218.04762
$ws.Range("I2").Value = 94
$ws.Range("J2").Value = 294.3846
$ws.Range("K2").Value = 94
$ws.Range("L2").Value = 294.3846
$ws.Range("M2").Value = 19
$ws.Range("N2").Value = -520.3846
$ws.Range("H5").Value = 1900
$ws.Range("I5").Value = 1900
$ws.Range("K5").Value = 1900
$ws.Range("M5").Value = -1788
$ws.Range("H42").Value = 31179
$ws.Range("J42").Value = 49995
$ws.Range("L42").Value = 49995
$ws.Range("N42").Value = -50965
$ws.Range("H49").Value = 25000.334
$ws.Range("J49").Value = 25000.334
$ws.Range("L49").Value = 25000.334
$ws.Range("N49").Value = -25368.334
$ws.Range("H55").Value = 34343.668
$ws.Range("I55").Value = 34998
$ws.Range("J55").Value = 34016.5
$ws.Range("K55").Value = 34998
$ws.Range("L55").Value = 34016.5
$ws.Range("M55").Value = -34671
$ws.Range("N55").Value = -34670.5
$ws.Range("H75").Value = 62218.25
$ws.Range("J75").Value = 62218.25
$ws.Range("L75").Value = 62218.25
$ws.Range("N75").Value = -63966.25
$ws.Range("H78").Value = 62218.25
$ws.Range("J78").Value = 62218.25
$ws.Range("L78").Value = 186654.75
$ws.Range("N78").Value = -195390.75
$ws.Range("H97").Value = 870.93335
$ws.Range("I97").Value = 557
$ws.Range("K97").Value = 557
$ws.Range("M97").Value = -61
$ws.Range("H102").Value = 4391
$ws.Range("I102").Value = 4252.1924
$ws.Range("J102").Value = 8000
$ws.Range("K102").Value = 4252.1924
$ws.Range("L102").Value = 8000
$ws.Range("M102").Value = -2630.1924
$ws.Range("N102").Value = -11244
$ws.Range("H115").Value = 31179
$ws.Range("J115").Value = 49995
$ws.Range("L115").Value = 49995
$ws.Range("N115").Value = -52345
$ws.Range("H122").Value = 14629.223
$ws.Range("I122").Value = 19194.75
$ws.Range("J122").Value = 5498.1665
$ws.Range("K122").Value = 57584.25
$ws.Range("L122").Value = 16494.4995
$ws.Range("M122").Value = -55134.25
$ws.Range("N122").Value = -21394.4995
$ws.Range("H123").Value = 46665.668
$ws.Range("J123").Value = 46665.668
$ws.Range("L123").Value = 46665.668
$ws.Range("N123").Value = -51565.668
$ws.Range("H132").Value = 4341.6665
$ws.Range("I132").Value = 4610
$ws.Range("K132").Value = 13830
$ws.Range("M132").Value = -11300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 10704
$ws.Range("I32").Value = 10704
$ws.Range("K32").Value = 10704
$ws.Range("M32").Value = -10387
$ws.Range("H40").Value = 6034.0938
$ws.Range("I40").Value = 4564.4707
$ws.Range("J40").Value = 7699.6665
$ws.Range("K40").Value = 4564.4707
$ws.Range("L40").Value = 7699.6665
$ws.Range("M40").Value = -4428.4707
$ws.Range("N40").Value = -7971.6665
$ws.Range("H93").Value = 4939.0527
$ws.Range("I93").Value = 2384.9
$ws.Range("K93").Value = 2384.9
$ws.Range("M93").Value = -1136.9
$ws.Range("H100").Value = 8714.286
$ws.Range("J100").Value = 9600
$ws.Range("L100").Value = 9600
$ws.Range("N100").Value = -10682
$ws.Range("H122").Value = 5698.9165
$ws.Range("I122").Value = 4923.375
$ws.Range("J122").Value = 7250
$ws.Range("K122").Value = 14770.125
$ws.Range("L122").Value = 21750
$ws.Range("M122").Value = -12320.125
$ws.Range("N122").Value = -26650
$ws.Range("H132").Value = 3749.5
$ws.Range("I132").Value = 4284.857
$ws.Range("K132").Value = 12854.571
$ws.Range("M132").Value = -10324.571
$ws.Range("H136").Value = 4369.5186
$ws.Range("J136").Value = 4012.2144
$ws.Range("L136").Value = 12036.6432
$ws.Range("N136").Value = -17136.6432
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 66668332
$ws.Range("I5").Value = 5000
$ws.Range("K5").Value = 5000
$ws.Range("M5").Value = -4888
$ws.Range("H81").Value = 3999.8
$ws.Range("I81").Value = 3999.8
$ws.Range("K81").Value = 7999.6
$ws.Range("M81").Value = -6938.6
$ws.Range("H84").Value = 3999.8
$ws.Range("I84").Value = 3999.8
$ws.Range("K84").Value = 39998
$ws.Range("M84").Value = -34694
$ws.Range("H100").Value = 802.5
$ws.Range("I100").Value = 802.5
$ws.Range("K100").Value = 1605
$ws.Range("M100").Value = -1064
$ws.Range("H107").Value = 2484.0527
$ws.Range("I107").Value = 1949.8125
$ws.Range("K107").Value = 5849.4375
$ws.Range("M107").Value = -3929.4375
$ws.Range("H119").Value = 110000
$ws.Range("J119").Value = 110000
$ws.Range("L119").Value = 110000
$ws.Range("N119").Value = -119676
$ws.Range("H122").Value = 5406.4346
$ws.Range("I122").Value = 2334.5386
$ws.Range("J122").Value = 9399.9
$ws.Range("K122").Value = 7003.6158
$ws.Range("L122").Value = 28199.7
$ws.Range("M122").Value = -4553.6158
$ws.Range("N122").Value = -33099.7
$ws.Range("H126").Value = 1344
$ws.Range("I126").Value = 1344
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4032
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1562
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1225.55
$ws.Range("I132").Value = 1108.4595
$ws.Range("K132").Value = 3325.3785
$ws.Range("M132").Value = -795.3784999999998
$ws.Range("H136").Value = 1583
$ws.Range("I136").Value = 1499.3334
$ws.Range("J136").Value = 1666.6666
$ws.Range("K136").Value = 4498.0002
$ws.Range("L136").Value = 4999.9998
$ws.Range("M136").Value = -1948.0002
$ws.Range("N136").Value = -10099.9998

Write-Host "Applied all changes"